$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 13.32779766666667
$ws.Range("H2").Value = 39.983393
$ws.Range("I2").Value = 0.1697233513642653
$ws.Range("J2").Value = 0.1697233513642653
$ws.Range("M2").Value = 0.04133166666666666
$ws.Range("Q2").Value = 0.5508600905594444
$ws.Range("R2").Value = 4.957740815035
$ws.Range("S2").Value = 0.1697233513642653
$ws.Range("T2").Value = 0.1697233513642653

# Row 3
$ws.Range("I3").Value = 0.5514955210569645
$ws.Range("J3").Value = 0.5514955210569645
$ws.Range("M3").Value = 0.04133166666666666
$ws.Range("S3").Value = 0.5514955210569645
$ws.Range("T3").Value = 0.5514955210569645

# Row 4
$ws.Range("H4").Value = 65.67520200000001
$ws.Range("I4").Value = 0.27878112757877
$ws.Range("J4").Value = 0.27878112757877
$ws.Range("M4").Value = 0.04133166666666666
$ws.Range("S4").Value = 0.27878112757877
$ws.Range("T4").Value = 0.27878112757877
